# Updates the "Clasificación" standings table with the results of another
# round of matches (PJ/games-played bumped from 5 to 6 for every player,
# plus the knock-on changes to wins/losses/goal-difference/points).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - David (Pos 1)
$ws.Range("C2").Value = 6.0
$ws.Range("F2").Value = 2.0
$ws.Range("H2").Value = 2.0
$ws.Range("I2").Value = 1.0

# Row 3 - Pedro (Pos 2)
$ws.Range("C3").Value = 6.0
$ws.Range("D3").Value = 2.0
$ws.Range("G3").Value = 3.0
$ws.Range("I3").Value = -4.0
$ws.Range("J3").Value = 2.0
$ws.Range("M3").Value = 6.0

# Row 4 - Adonay (Pos 3)
$ws.Range("C4").Value = 6.0
$ws.Range("D4").Value = 6.0
$ws.Range("G4").Value = 7.0
$ws.Range("I4").Value = 7.0
$ws.Range("J4").Value = 5.0
$ws.Range("M4").Value = 19.0

# Row 5 - Richard (Pos 4)
$ws.Range("C5").Value = 6.0
$ws.Range("F5").Value = 5.0
$ws.Range("H5").Value = 5.0
$ws.Range("I5").Value = -5.0

# Row 6 - Iván (Pos 5)
$ws.Range("C6").Value = 6.0
$ws.Range("F6").Value = 3.0
$ws.Range("H6").Value = 4.0
$ws.Range("I6").Value = 3.0

# Row 7 - Nico (Pos 6)
$ws.Range("C7").Value = 6.0
$ws.Range("F7").Value = 3.0
$ws.Range("H7").Value = 3.0
$ws.Range("I7").Value = -1.0

# Row 8 - Nicolás (Pos 7)
$ws.Range("C8").Value = 6.0
$ws.Range("D8").Value = 2.0
$ws.Range("G8").Value = 2.0
$ws.Range("I8").Value = -4.0
$ws.Range("J8").Value = 2.0
$ws.Range("M8").Value = 6.0

# Row 9 - Vicente (Pos 8)
$ws.Range("C9").Value = 6.0
$ws.Range("D9").Value = 4.0
$ws.Range("G9").Value = 4.0
$ws.Range("I9").Value = 2.0
$ws.Range("J9").Value = 4.0
$ws.Range("M9").Value = 9.0
